$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (B1:J1)
$ws.Cells.Item(1, 2).Value = "log_name"
$ws.Cells.Item(1, 3).Value = "C3"
$ws.Cells.Item(1, 4).Value = "C4"
$ws.Cells.Item(1, 5).Value = "C5"
$ws.Cells.Item(1, 6).Value = "N1"
$ws.Cells.Item(1, 7).Value = "C6"
$ws.Cells.Item(1, 8).Value = "C7"
$ws.Cells.Item(1, 9).Value = "C1"
$ws.Cells.Item(1, 10).Value = "C2"

# Apply header style (bold, centered, bordered) to the new I1:J1 header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-38 (columns B..J only; column A unchanged)
# Row 2
$ws.Cells.Item(2, 2).Value = "pyrd10_conf-1"
$ws.Cells.Item(2, 3).Value = "C10"
$ws.Cells.Item(2, 4).Value = "C5"
$ws.Cells.Item(2, 5).Value = "C6"
$ws.Cells.Item(2, 6).Value = "N7"
$ws.Cells.Item(2, 7).Value = "C8"
$ws.Cells.Item(2, 8).Value = "C9"
$ws.Cells.Item(2, 9).Value = "C4"
$ws.Cells.Item(2, 10).Value = "C5"
# Row 3
$ws.Cells.Item(3, 2).Value = "pyrd10_conf-10"
$ws.Cells.Item(3, 3).Value = "C10"
$ws.Cells.Item(3, 4).Value = "C5"
$ws.Cells.Item(3, 5).Value = "C6"
$ws.Cells.Item(3, 6).Value = "N7"
$ws.Cells.Item(3, 7).Value = "C8"
$ws.Cells.Item(3, 8).Value = "C9"
$ws.Cells.Item(3, 9).Value = "C4"
$ws.Cells.Item(3, 10).Value = "C5"
# Row 4
$ws.Cells.Item(4, 2).Value = "pyrd10_conf-11"
$ws.Cells.Item(4, 3).Value = "C10"
$ws.Cells.Item(4, 4).Value = "C5"
$ws.Cells.Item(4, 5).Value = "C6"
$ws.Cells.Item(4, 6).Value = "N7"
$ws.Cells.Item(4, 7).Value = "C8"
$ws.Cells.Item(4, 8).Value = "C9"
$ws.Cells.Item(4, 9).Value = "C4"
$ws.Cells.Item(4, 10).Value = "C5"
# Row 5
$ws.Cells.Item(5, 2).Value = "pyrd10_conf-12"
$ws.Cells.Item(5, 3).Value = "C10"
$ws.Cells.Item(5, 4).Value = "C9"
$ws.Cells.Item(5, 5).Value = "C8"
$ws.Cells.Item(5, 6).Value = "N7"
$ws.Cells.Item(5, 7).Value = "C6"
$ws.Cells.Item(5, 8).Value = "C5"
$ws.Cells.Item(5, 9).Value = "C4"
$ws.Cells.Item(5, 10).Value = "C5"
# Row 6
$ws.Cells.Item(6, 2).Value = "pyrd10_conf-13"
$ws.Cells.Item(6, 3).Value = "C10"
$ws.Cells.Item(6, 4).Value = "C5"
$ws.Cells.Item(6, 5).Value = "C6"
$ws.Cells.Item(6, 6).Value = "N7"
$ws.Cells.Item(6, 7).Value = "C8"
$ws.Cells.Item(6, 8).Value = "C9"
$ws.Cells.Item(6, 9).Value = "C4"
$ws.Cells.Item(6, 10).Value = "C5"
# Row 7
$ws.Cells.Item(7, 2).Value = "pyrd10_conf-2"
$ws.Cells.Item(7, 3).Value = "C10"
$ws.Cells.Item(7, 4).Value = "C5"
$ws.Cells.Item(7, 5).Value = "C6"
$ws.Cells.Item(7, 6).Value = "N7"
$ws.Cells.Item(7, 7).Value = "C8"
$ws.Cells.Item(7, 8).Value = "C9"
$ws.Cells.Item(7, 9).Value = "C4"
$ws.Cells.Item(7, 10).Value = "C5"
# Row 8
$ws.Cells.Item(8, 2).Value = "pyrd10_conf-3"
$ws.Cells.Item(8, 3).Value = "C10"
$ws.Cells.Item(8, 4).Value = "C5"
$ws.Cells.Item(8, 5).Value = "C6"
$ws.Cells.Item(8, 6).Value = "N7"
$ws.Cells.Item(8, 7).Value = "C8"
$ws.Cells.Item(8, 8).Value = "C9"
$ws.Cells.Item(8, 9).Value = "C4"
$ws.Cells.Item(8, 10).Value = "C5"
# Row 9
$ws.Cells.Item(9, 2).Value = "pyrd10_conf-4"
$ws.Cells.Item(9, 3).Value = "C10"
$ws.Cells.Item(9, 4).Value = "C5"
$ws.Cells.Item(9, 5).Value = "C6"
$ws.Cells.Item(9, 6).Value = "N7"
$ws.Cells.Item(9, 7).Value = "C8"
$ws.Cells.Item(9, 8).Value = "C9"
$ws.Cells.Item(9, 9).Value = "C4"
$ws.Cells.Item(9, 10).Value = "C5"
# Row 10
$ws.Cells.Item(10, 2).Value = "pyrd10_conf-5"
$ws.Cells.Item(10, 3).Value = "C10"
$ws.Cells.Item(10, 4).Value = "C5"
$ws.Cells.Item(10, 5).Value = "C6"
$ws.Cells.Item(10, 6).Value = "N7"
$ws.Cells.Item(10, 7).Value = "C8"
$ws.Cells.Item(10, 8).Value = "C9"
$ws.Cells.Item(10, 9).Value = "C4"
$ws.Cells.Item(10, 10).Value = "C5"
# Row 11
$ws.Cells.Item(11, 2).Value = "pyrd10_conf-6"
$ws.Cells.Item(11, 3).Value = "C10"
$ws.Cells.Item(11, 4).Value = "C5"
$ws.Cells.Item(11, 5).Value = "C6"
$ws.Cells.Item(11, 6).Value = "N7"
$ws.Cells.Item(11, 7).Value = "C8"
$ws.Cells.Item(11, 8).Value = "C9"
$ws.Cells.Item(11, 9).Value = "C4"
$ws.Cells.Item(11, 10).Value = "C5"
# Row 12
$ws.Cells.Item(12, 2).Value = "pyrd10_conf-7"
$ws.Cells.Item(12, 3).Value = "C10"
$ws.Cells.Item(12, 4).Value = "C5"
$ws.Cells.Item(12, 5).Value = "C6"
$ws.Cells.Item(12, 6).Value = "N7"
$ws.Cells.Item(12, 7).Value = "C8"
$ws.Cells.Item(12, 8).Value = "C9"
$ws.Cells.Item(12, 9).Value = "C4"
$ws.Cells.Item(12, 10).Value = "C5"
# Row 13
$ws.Cells.Item(13, 2).Value = "pyrd10_conf-8"
$ws.Cells.Item(13, 3).Value = "C10"
$ws.Cells.Item(13, 4).Value = "C5"
$ws.Cells.Item(13, 5).Value = "C6"
$ws.Cells.Item(13, 6).Value = "N7"
$ws.Cells.Item(13, 7).Value = "C8"
$ws.Cells.Item(13, 8).Value = "C9"
$ws.Cells.Item(13, 9).Value = "C4"
$ws.Cells.Item(13, 10).Value = "C5"
# Row 14
$ws.Cells.Item(14, 2).Value = "pyrd10_conf-9"
$ws.Cells.Item(14, 3).Value = "C10"
$ws.Cells.Item(14, 4).Value = "C9"
$ws.Cells.Item(14, 5).Value = "C8"
$ws.Cells.Item(14, 6).Value = "N7"
$ws.Cells.Item(14, 7).Value = "C6"
$ws.Cells.Item(14, 8).Value = "C5"
$ws.Cells.Item(14, 9).Value = "C4"
$ws.Cells.Item(14, 10).Value = "C5"
# Row 15
$ws.Cells.Item(15, 2).Value = "pyrd11_conf-1"
$ws.Cells.Item(15, 3).Value = "C2"
$ws.Cells.Item(15, 4).Value = "C12"
$ws.Cells.Item(15, 5).Value = "C11"
$ws.Cells.Item(15, 6).Value = "N10"
$ws.Cells.Item(15, 7).Value = "C5"
$ws.Cells.Item(15, 8).Value = "C3"
$ws.Cells.Item(15, 9).Value = "C4"
$ws.Cells.Item(15, 10).Value = "C3"
# Row 16
$ws.Cells.Item(16, 2).Value = "pyrd11_conf-2"
$ws.Cells.Item(16, 3).Value = "C2"
$ws.Cells.Item(16, 4).Value = "C3"
$ws.Cells.Item(16, 5).Value = "C5"
$ws.Cells.Item(16, 6).Value = "N10"
$ws.Cells.Item(16, 7).Value = "C11"
$ws.Cells.Item(16, 8).Value = "C12"
$ws.Cells.Item(16, 9).Value = "C4"
$ws.Cells.Item(16, 10).Value = "C3"
# Row 17
$ws.Cells.Item(17, 2).Value = "pyrd12_conf-1"
$ws.Cells.Item(17, 3).Value = "C3"
$ws.Cells.Item(17, 4).Value = "C13"
$ws.Cells.Item(17, 5).Value = "C12"
$ws.Cells.Item(17, 6).Value = "N11"
$ws.Cells.Item(17, 7).Value = "C10"
$ws.Cells.Item(17, 8).Value = "C4"
$ws.Cells.Item(17, 9).Value = "C5"
$ws.Cells.Item(17, 10).Value = "C4"
# Row 18
$ws.Cells.Item(18, 2).Value = "pyrd13_conf-1"
$ws.Cells.Item(18, 3).Value = "C5"
$ws.Cells.Item(18, 4).Value = "C4"
$ws.Cells.Item(18, 5).Value = "C3"
$ws.Cells.Item(18, 6).Value = "N12"
$ws.Cells.Item(18, 7).Value = "C11"
$ws.Cells.Item(18, 8).Value = "C6"
$ws.Cells.Item(18, 9).Value = "C7"
$ws.Cells.Item(18, 10).Value = "C6"
# Row 19
$ws.Cells.Item(19, 2).Value = "pyrd1_conf-1"
$ws.Cells.Item(19, 3).Value = "C3"
$ws.Cells.Item(19, 4).Value = "C2"
$ws.Cells.Item(19, 5).Value = "C11"
$ws.Cells.Item(19, 6).Value = "N10"
$ws.Cells.Item(19, 7).Value = "C9"
$ws.Cells.Item(19, 8).Value = "C4"
$ws.Cells.Item(19, 9).Value = "C8"
$ws.Cells.Item(19, 10).Value = "C3"
# Row 20
$ws.Cells.Item(20, 2).Value = "pyrd2_conf-1"
$ws.Cells.Item(20, 3).Value = "C6"
$ws.Cells.Item(20, 4).Value = "C5"
$ws.Cells.Item(20, 5).Value = "C4"
$ws.Cells.Item(20, 6).Value = "N3"
$ws.Cells.Item(20, 7).Value = "C2"
$ws.Cells.Item(20, 8).Value = "C7"
$ws.Cells.Item(20, 9).Value = "C8"
$ws.Cells.Item(20, 10).Value = "C7"
# Row 21
$ws.Cells.Item(21, 2).Value = "pyrd3_conf-1"
$ws.Cells.Item(21, 3).Value = "C2"
$ws.Cells.Item(21, 4).Value = "C8"
$ws.Cells.Item(21, 5).Value = "C7"
$ws.Cells.Item(21, 6).Value = "N6"
$ws.Cells.Item(21, 7).Value = "C4"
$ws.Cells.Item(21, 8).Value = "C3"
$ws.Cells.Item(21, 9).Value = "C5"
$ws.Cells.Item(21, 10).Value = "C4"
# Row 22
$ws.Cells.Item(22, 2).Value = "pyrd4_conf-1"
$ws.Cells.Item(22, 3).Value = "C8"
$ws.Cells.Item(22, 4).Value = "C7"
$ws.Cells.Item(22, 5).Value = "C6"
$ws.Cells.Item(22, 6).Value = "N5"
$ws.Cells.Item(22, 7).Value = "C4"
$ws.Cells.Item(22, 8).Value = "C12"
$ws.Cells.Item(22, 9).Value = "C9"
$ws.Cells.Item(22, 10).Value = "C8"
# Row 23
$ws.Cells.Item(23, 2).Value = "pyrd4_conf-2"
$ws.Cells.Item(23, 3).Value = "C8"
$ws.Cells.Item(23, 4).Value = "C7"
$ws.Cells.Item(23, 5).Value = "C6"
$ws.Cells.Item(23, 6).Value = "N5"
$ws.Cells.Item(23, 7).Value = "C4"
$ws.Cells.Item(23, 8).Value = "C12"
$ws.Cells.Item(23, 9).Value = "C9"
$ws.Cells.Item(23, 10).Value = "C8"
# Row 24
$ws.Cells.Item(24, 2).Value = "pyrd4_conf-3"
$ws.Cells.Item(24, 3).Value = "C8"
$ws.Cells.Item(24, 4).Value = "C7"
$ws.Cells.Item(24, 5).Value = "C6"
$ws.Cells.Item(24, 6).Value = "N5"
$ws.Cells.Item(24, 7).Value = "C4"
$ws.Cells.Item(24, 8).Value = "C12"
$ws.Cells.Item(24, 9).Value = "C9"
$ws.Cells.Item(24, 10).Value = "C8"
# Row 25
$ws.Cells.Item(25, 2).Value = "pyrd4_conf-4"
$ws.Cells.Item(25, 3).Value = "C8"
$ws.Cells.Item(25, 4).Value = "C7"
$ws.Cells.Item(25, 5).Value = "C6"
$ws.Cells.Item(25, 6).Value = "N5"
$ws.Cells.Item(25, 7).Value = "C4"
$ws.Cells.Item(25, 8).Value = "C12"
$ws.Cells.Item(25, 9).Value = "C9"
$ws.Cells.Item(25, 10).Value = "C8"
# Row 26
$ws.Cells.Item(26, 2).Value = "pyrd4_conf-5"
$ws.Cells.Item(26, 3).Value = "C8"
$ws.Cells.Item(26, 4).Value = "C7"
$ws.Cells.Item(26, 5).Value = "C6"
$ws.Cells.Item(26, 6).Value = "N5"
$ws.Cells.Item(26, 7).Value = "C4"
$ws.Cells.Item(26, 8).Value = "C12"
$ws.Cells.Item(26, 9).Value = "C9"
$ws.Cells.Item(26, 10).Value = "C8"
# Row 27
$ws.Cells.Item(27, 2).Value = "pyrd4_conf-6"
$ws.Cells.Item(27, 3).Value = "C8"
$ws.Cells.Item(27, 4).Value = "C7"
$ws.Cells.Item(27, 5).Value = "C6"
$ws.Cells.Item(27, 6).Value = "N5"
$ws.Cells.Item(27, 7).Value = "C4"
$ws.Cells.Item(27, 8).Value = "C12"
$ws.Cells.Item(27, 9).Value = "C9"
$ws.Cells.Item(27, 10).Value = "C8"
# Row 28
$ws.Cells.Item(28, 2).Value = "pyrd4_conf-7"
$ws.Cells.Item(28, 3).Value = "C8"
$ws.Cells.Item(28, 4).Value = "C7"
$ws.Cells.Item(28, 5).Value = "C6"
$ws.Cells.Item(28, 6).Value = "N5"
$ws.Cells.Item(28, 7).Value = "C4"
$ws.Cells.Item(28, 8).Value = "C12"
$ws.Cells.Item(28, 9).Value = "C9"
$ws.Cells.Item(28, 10).Value = "C8"
# Row 29
$ws.Cells.Item(29, 2).Value = "pyrd4_conf-8"
$ws.Cells.Item(29, 3).Value = "C8"
$ws.Cells.Item(29, 4).Value = "C7"
$ws.Cells.Item(29, 5).Value = "C6"
$ws.Cells.Item(29, 6).Value = "N5"
$ws.Cells.Item(29, 7).Value = "C4"
$ws.Cells.Item(29, 8).Value = "C12"
$ws.Cells.Item(29, 9).Value = "C9"
$ws.Cells.Item(29, 10).Value = "C8"
# Row 30
$ws.Cells.Item(30, 2).Value = "pyrd5_conf-1"
$ws.Cells.Item(30, 3).Value = "C10"
$ws.Cells.Item(30, 4).Value = "C9"
$ws.Cells.Item(30, 5).Value = "C4"
$ws.Cells.Item(30, 6).Value = "N3"
$ws.Cells.Item(30, 7).Value = "C2"
$ws.Cells.Item(30, 8).Value = "C12"
$ws.Cells.Item(30, 9).Value = "C11"
$ws.Cells.Item(30, 10).Value = "C10"
# Row 31
$ws.Cells.Item(31, 2).Value = "pyrd6_conf-1"
$ws.Cells.Item(31, 3).Value = "C4"
$ws.Cells.Item(31, 4).Value = "C6"
$ws.Cells.Item(31, 5).Value = "C7"
$ws.Cells.Item(31, 6).Value = "N16"
$ws.Cells.Item(31, 7).Value = "C2"
$ws.Cells.Item(31, 8).Value = "C3"
$ws.Cells.Item(31, 9).Value = "C5"
$ws.Cells.Item(31, 10).Value = "C4"
# Row 32
$ws.Cells.Item(32, 2).Value = "pyrd6_conf-2"
$ws.Cells.Item(32, 3).Value = "C4"
$ws.Cells.Item(32, 4).Value = "C3"
$ws.Cells.Item(32, 5).Value = "C2"
$ws.Cells.Item(32, 6).Value = "N16"
$ws.Cells.Item(32, 7).Value = "C7"
$ws.Cells.Item(32, 8).Value = "C6"
$ws.Cells.Item(32, 9).Value = "C5"
$ws.Cells.Item(32, 10).Value = "C4"
# Row 33
$ws.Cells.Item(33, 2).Value = "pyrd6_conf-3"
$ws.Cells.Item(33, 3).Value = "C4"
$ws.Cells.Item(33, 4).Value = "C3"
$ws.Cells.Item(33, 5).Value = "C2"
$ws.Cells.Item(33, 6).Value = "N16"
$ws.Cells.Item(33, 7).Value = "C7"
$ws.Cells.Item(33, 8).Value = "C6"
$ws.Cells.Item(33, 9).Value = "C5"
$ws.Cells.Item(33, 10).Value = "C4"
# Row 34
$ws.Cells.Item(34, 2).Value = "pyrd7_conf-1"
$ws.Cells.Item(34, 3).Value = "C19"
$ws.Cells.Item(34, 4).Value = "C13"
$ws.Cells.Item(34, 5).Value = "C14"
$ws.Cells.Item(34, 6).Value = "N15"
$ws.Cells.Item(34, 7).Value = "C16"
$ws.Cells.Item(34, 8).Value = "C17"
$ws.Cells.Item(34, 9).Value = "C18"
$ws.Cells.Item(34, 10).Value = "C17"
# Row 35
$ws.Cells.Item(35, 2).Value = "pyrd7_conf-2"
$ws.Cells.Item(35, 3).Value = "C19"
$ws.Cells.Item(35, 4).Value = "C13"
$ws.Cells.Item(35, 5).Value = "C14"
$ws.Cells.Item(35, 6).Value = "N15"
$ws.Cells.Item(35, 7).Value = "C16"
$ws.Cells.Item(35, 8).Value = "C17"
$ws.Cells.Item(35, 9).Value = "C18"
$ws.Cells.Item(35, 10).Value = "C17"
# Row 36
$ws.Cells.Item(36, 2).Value = "pyrd8_conf-1"
$ws.Cells.Item(36, 3).Value = "C11"
$ws.Cells.Item(36, 4).Value = "C10"
$ws.Cells.Item(36, 5).Value = "C4"
$ws.Cells.Item(36, 6).Value = "N3"
$ws.Cells.Item(36, 7).Value = "C2"
$ws.Cells.Item(36, 8).Value = "C13"
$ws.Cells.Item(36, 9).Value = "C12"
$ws.Cells.Item(36, 10).Value = "C11"
# Row 37
$ws.Cells.Item(37, 2).Value = "pyrd8_conf-2"
$ws.Cells.Item(37, 3).Value = "C11"
$ws.Cells.Item(37, 4).Value = "C13"
$ws.Cells.Item(37, 5).Value = "C2"
$ws.Cells.Item(37, 6).Value = "N3"
$ws.Cells.Item(37, 7).Value = "C4"
$ws.Cells.Item(37, 8).Value = "C10"
$ws.Cells.Item(37, 9).Value = "C12"
$ws.Cells.Item(37, 10).Value = "C11"
# Row 38
$ws.Cells.Item(38, 2).Value = "pyrd9_conf-1"
$ws.Cells.Item(38, 3).Value = "C7"
$ws.Cells.Item(38, 4).Value = "C6"
$ws.Cells.Item(38, 5).Value = "C4"
$ws.Cells.Item(38, 6).Value = "N3"
$ws.Cells.Item(38, 7).Value = "C2"
$ws.Cells.Item(38, 8).Value = "C8"
$ws.Cells.Item(38, 9).Value = "C1"
$ws.Cells.Item(38, 10).Value = "C2"
